$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values
$ws.Range("C2").Value = 360
$ws.Range("C3").Value = 82.7777
$ws.Range("C5").Value = 127.873229

# Update selection on the sheet (active cell moves to D8)
$ws.Range("D8").Select()

# Adjust the workbook window tab ratio (stored as a fraction, 0.655 == 65.5%)
$wb.Windows.Item(1).TabRatio = 0.655
